$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.711.67'
$ws.Range('E2').Value = '  -1.78%  '
$ws.Range('D3').Value = '3.382.49'
$ws.Range('E3').Value = '  -2.21%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.12'
$ws.Range('E5').Value = '  -1.59%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.51'
$ws.Range('E6').Value = '  -4.73%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '3.383.16'
$ws.Range('E8').Value = '  -2.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.473'
$ws.Range('E9').Value = '  -0.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.50'
$ws.Range('E11').Value = '  -0.95%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.393'
$ws.Range('E12').Value = '  +0.38%  '
$ws.Range('D13').Value = '3.961.03'
$ws.Range('E13').Value = '  -2.21%  '
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.124'
$ws.Range('E14').Value = '  +0.74%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '28.00'
$ws.Range('E15').Value = '  -0.83%  '
$ws.Range('E16').Value = '  -2.74%  '
$ws.Range('D17').Value = '3.388.08'
$ws.Range('E17').Value = '  -1.64%  '
$ws.Range('D18').Value = '60.832.22'
$ws.Range('E18').Value = '  -1.66%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.29'
$ws.Range('E19').Value = '  -0.91%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.09'
$ws.Range('E20').Value = '  -1.57%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.94'
$ws.Range('E21').Value = '  -5.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '388.54'
$ws.Range('E22').Value = '  +0.78%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.560'
$ws.Range('E23').Value = '  -1.80%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.24'
$ws.Range('E24').Value = '  +0.87%  '
$ws.Range('E25').Value = '  -0.24%  '
$ws.Range('E26').Value = '  -3.92%  '
$ws.Range('D27').Value = '3.524.55'
$ws.Range('E27').Value = '  -2.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.178'
$ws.Range('E28').Value = '  -1.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('E30').Value = '  -5.76%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.07'
$ws.Range('E31').Value = '  -2.13%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.45'
$ws.Range('E32').Value = '  -5.01%  '
$ws.Range('E33').Value = '  -1.55%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.76'
$ws.Range('E35').Value = '  -0.82%  '
$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.91'
$ws.Range('E36').Value = '  -2.14%  '
$ws.Range('B37').Value = 'RenzoRestakedETH'
$ws.Range('C37').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D37').Value = '3.408.71'
$ws.Range('E37').Value = '  -2.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '166.76'
$ws.Range('E38').Value = '  +0.14%  '
$ws.Range('E39').Value = '  -3.34%  '
$ws.Range('E40').Value = '  -3.15%  '
$ws.Range('E41').Value = '  -1.52%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '26.90'
$ws.Range('E42').Value = '  +3.45%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.783'
$ws.Range('E43').Value = '  -1.76%  '
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('E45').Value = '  -0.71%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '41.68'
$ws.Range('E46').Value = '  -1.59%  '
$ws.Range('E47').Value = '  -2.61%  '
$ws.Range('D48').Value = '2.547.55'
$ws.Range('E48').Value = '  -2.35%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.12'
$ws.Range('E49').Value = '  -4.39%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.82'
$ws.Range('E50').Value = '  -1.96%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '22.86'
$ws.Range('E51').Value = '  -2.11%  '
